# Scrum Master: Files Update
# Add a new task row to the "week1" sheet and move the selection down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week1")

# New task text in B12 (continues the list that ends at B11)
$ws.Range("B12").Value = "Pesquisar como se joga o jogo"

# B13 is left blank but carries the same "underline placeholder" style as B15
$ws.Range("B13").Value = ""
$ws.Range("B13").Font.Underline = 2

# Move/confirm the active selection onto B13, like typing Enter after B12
$ws.Activate()
$ws.Range("B13").Select()
